$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") should carry the same style as
# the other header cells (bold font, thin border, centered). Copy the
# formatting from the existing H1 header cell, then overwrite the values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-15 for the two new columns I (I0) and J (IF).
$values = @{
    2  = @(8, 9)
    3  = @(1, 3)
    4  = @(2, 4)
    5  = @(5, 9)
    6  = @(7, 7)
    7  = @(1, 3)
    8  = @(1, 5)
    9  = @(2, 6)
    10 = @(4, 7)
    11 = @(1, 2)
    12 = @(1, 4)
    13 = @(1, 3)
    14 = @(1, 2)
    15 = @(1, 1)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
